$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Adform")

# --- Update the row-2 field-requirement descriptions to mention the new
#     "Delete" capability (added alongside the existing Add/Edit rules). ---
$ws.Range("A2").Value = "Add: Not required`nEdit: Required`nDelete: Required"
$ws.Range("B2").Value = "Add: Required`nEdit: Required`nDelete: Not Required"
$ws.Range("C2").Value = "Add: Required`nEdit: Required`nDelete: Not Required"
$ws.Range("D2").Value = "Add: Required`nEdit: Required`nDelete: Not Required`nValues: global, apac`nDefaulted to global"
$ws.Range("E2").Value = "Add: Required`nEdit: Required`nDelete: Not Required"
$ws.Range("F2").Value = "Add: Required`nEdit: Required`nDelete: Not Required"
$ws.Range("G2").Value = "Add: Required`nEdit: Required`nDelete: Not Required"

# --- Row 2 is now taller to fit the extra "Delete" line. ---
$ws.Rows("2").RowHeight = 82

# --- Widen the columns that grew to comfortably fit the longer text. ---
$ws.Columns("B").ColumnWidth = 17.833333333333332
$ws.Columns("C").ColumnWidth = 18
$ws.Columns("E").ColumnWidth = 18.833333333333332
$ws.Range("F:G").ColumnWidth = 18.166666666666668

# --- Make Adform the active sheet / selection (it was TTD before). ---
$ws.Activate()
$ws.Range("A2:I2").Select()
